$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (source data from original row 2)
$ws.Cells.Item(2, 3).Value = 46079

# Row 3 (source data from original row 3)
$ws.Cells.Item(3, 3).Value = 46079

# Row 4 (source data from original row 4)
$ws.Cells.Item(4, 3).Value = 46079

# Row 5 (source data from original row 10)
$ws.Cells.Item(5, 1).Value = 'A 34341-2024'
$ws.Cells.Item(5, 2).Value = 45525
$ws.Cells.Item(5, 3).Value = 46079
$ws.Cells.Item(5, 6).Value = 'Övriga Aktiebolag'
$ws.Cells.Item(5, 7).Value = 14.4
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = 0
$ws.Cells.Item(5, 15).Value = 1
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = 'Desmeknopp'
$ws.Cells.Item(5, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 34341-2024 artfynd.xlsx", "A 34341-2024")'
$ws.Cells.Item(5, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 34341-2024 karta.png", "A 34341-2024")'
$ws.Cells.Item(5, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 34341-2024 FSC-klagomål.docx", "A 34341-2024")'
$ws.Cells.Item(5, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 34341-2024 FSC-klagomål mail.docx", "A 34341-2024")'
$ws.Cells.Item(5, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 34341-2024 tillsynsbegäran.docx", "A 34341-2024")'
$ws.Cells.Item(5, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 34341-2024 tillsynsbegäran mail.docx", "A 34341-2024")'

# Row 6 (source data from original row 9)
$ws.Cells.Item(6, 1).Value = 'A 31213-2023'
$ws.Cells.Item(6, 2).Value = 45113
$ws.Cells.Item(6, 3).Value = 46079
$ws.Cells.Item(6, 6).ClearContents()
$ws.Cells.Item(6, 7).Value = 6.5
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = 0
$ws.Cells.Item(6, 14).Value = 0
$ws.Cells.Item(6, 15).Value = 1
$ws.Cells.Item(6, 16).Value = 0
$ws.Cells.Item(6, 17).Value = 1
$ws.Cells.Item(6, 18).Value = 'Skogsveronika'
$ws.Cells.Item(6, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 31213-2023 artfynd.xlsx", "A 31213-2023")'
$ws.Cells.Item(6, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 31213-2023 karta.png", "A 31213-2023")'
$ws.Cells.Item(6, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 31213-2023 FSC-klagomål.docx", "A 31213-2023")'
$ws.Cells.Item(6, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 31213-2023 FSC-klagomål mail.docx", "A 31213-2023")'
$ws.Cells.Item(6, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 31213-2023 tillsynsbegäran.docx", "A 31213-2023")'
$ws.Cells.Item(6, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 31213-2023 tillsynsbegäran mail.docx", "A 31213-2023")'

# Row 7 (source data from original row 6)
$ws.Cells.Item(7, 1).Value = 'A 13766-2023'
$ws.Cells.Item(7, 2).Value = 45007
$ws.Cells.Item(7, 3).Value = 46079
$ws.Cells.Item(7, 6).ClearContents()
$ws.Cells.Item(7, 7).Value = 0.9
$ws.Cells.Item(7, 8).Value = 1
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = 0
$ws.Cells.Item(7, 14).Value = 0
$ws.Cells.Item(7, 15).Value = 0
$ws.Cells.Item(7, 16).Value = 0
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = 'Större vattensalamander'
$ws.Cells.Item(7, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 13766-2023 artfynd.xlsx", "A 13766-2023")'
$ws.Cells.Item(7, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 13766-2023 karta.png", "A 13766-2023")'
$ws.Cells.Item(7, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 13766-2023 FSC-klagomål.docx", "A 13766-2023")'
$ws.Cells.Item(7, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 13766-2023 FSC-klagomål mail.docx", "A 13766-2023")'
$ws.Cells.Item(7, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 13766-2023 tillsynsbegäran.docx", "A 13766-2023")'
$ws.Cells.Item(7, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 13766-2023 tillsynsbegäran mail.docx", "A 13766-2023")'

# Row 8 (source data from original row 7)
$ws.Cells.Item(8, 1).Value = 'A 60891-2024'
$ws.Cells.Item(8, 2).Value = 45644
$ws.Cells.Item(8, 3).Value = 46079
$ws.Cells.Item(8, 6).ClearContents()
$ws.Cells.Item(8, 7).Value = 16.1
$ws.Cells.Item(8, 8).Value = 1
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = 0
$ws.Cells.Item(8, 14).Value = 0
$ws.Cells.Item(8, 15).Value = 0
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(8, 17).Value = 1
$ws.Cells.Item(8, 18).Value = 'Lövgroda'
$ws.Cells.Item(8, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 60891-2024 artfynd.xlsx", "A 60891-2024")'
$ws.Cells.Item(8, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 60891-2024 karta.png", "A 60891-2024")'
$ws.Cells.Item(8, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 60891-2024 FSC-klagomål.docx", "A 60891-2024")'
$ws.Cells.Item(8, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 60891-2024 FSC-klagomål mail.docx", "A 60891-2024")'
$ws.Cells.Item(8, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 60891-2024 tillsynsbegäran.docx", "A 60891-2024")'
$ws.Cells.Item(8, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 60891-2024 tillsynsbegäran mail.docx", "A 60891-2024")'

# Row 9 (source data from original row 5)
$ws.Cells.Item(9, 1).Value = 'A 61558-2023'
$ws.Cells.Item(9, 2).Value = 45265
$ws.Cells.Item(9, 3).Value = 46079
$ws.Cells.Item(9, 6).Value = 'Övriga statliga verk och myndigheter'
$ws.Cells.Item(9, 7).Value = 1.5
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0
$ws.Cells.Item(9, 14).Value = 0
$ws.Cells.Item(9, 15).Value = 1
$ws.Cells.Item(9, 16).Value = 1
$ws.Cells.Item(9, 17).Value = 1
$ws.Cells.Item(9, 18).Value = 'Ask'
$ws.Cells.Item(9, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 61558-2023 artfynd.xlsx", "A 61558-2023")'
$ws.Cells.Item(9, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 61558-2023 karta.png", "A 61558-2023")'
$ws.Cells.Item(9, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 61558-2023 FSC-klagomål.docx", "A 61558-2023")'
$ws.Cells.Item(9, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 61558-2023 FSC-klagomål mail.docx", "A 61558-2023")'
$ws.Cells.Item(9, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 61558-2023 tillsynsbegäran.docx", "A 61558-2023")'
$ws.Cells.Item(9, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 61558-2023 tillsynsbegäran mail.docx", "A 61558-2023")'

# Row 10 (source data from original row 8)
$ws.Cells.Item(10, 1).Value = 'A 49546-2025'
$ws.Cells.Item(10, 2).Value = 45939
$ws.Cells.Item(10, 3).Value = 46079
$ws.Cells.Item(10, 6).ClearContents()
$ws.Cells.Item(10, 7).Value = 4.4
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 1
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = 0
$ws.Cells.Item(10, 14).Value = 0
$ws.Cells.Item(10, 15).Value = 0
$ws.Cells.Item(10, 16).Value = 0
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = 'Igelkottsröksvamp'
$ws.Cells.Item(10, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 49546-2025 artfynd.xlsx", "A 49546-2025")'
$ws.Cells.Item(10, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 49546-2025 karta.png", "A 49546-2025")'
$ws.Cells.Item(10, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 49546-2025 FSC-klagomål.docx", "A 49546-2025")'
$ws.Cells.Item(10, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 49546-2025 FSC-klagomål mail.docx", "A 49546-2025")'
$ws.Cells.Item(10, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 49546-2025 tillsynsbegäran.docx", "A 49546-2025")'
$ws.Cells.Item(10, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 49546-2025 tillsynsbegäran mail.docx", "A 49546-2025")'

# Row 11 (source data from original row 11)
$ws.Cells.Item(11, 3).Value = 46079

# Row 12 (source data from original row 12)
$ws.Cells.Item(12, 3).Value = 46079

# Row 13 (source data from original row 13)
$ws.Cells.Item(13, 3).Value = 46079

# Row 14 (source data from original row 19)
$ws.Cells.Item(14, 1).Value = 'A 60803-2023'
$ws.Cells.Item(14, 2).Value = 45260
$ws.Cells.Item(14, 3).Value = 46079
$ws.Cells.Item(14, 6).ClearContents()
$ws.Cells.Item(14, 7).Value = 1.6
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = 0
$ws.Cells.Item(14, 14).Value = 0
$ws.Cells.Item(14, 15).Value = 0
$ws.Cells.Item(14, 16).Value = 0
$ws.Cells.Item(14, 17).Value = 0
$ws.Cells.Item(14, 18).ClearContents()
$ws.Cells.Item(14, 19).ClearContents()
$ws.Cells.Item(14, 20).ClearContents()
$ws.Cells.Item(14, 22).ClearContents()
$ws.Cells.Item(14, 23).ClearContents()
$ws.Cells.Item(14, 24).ClearContents()
$ws.Cells.Item(14, 25).ClearContents()

# Row 15 (source data from original row 17)
$ws.Cells.Item(15, 1).Value = 'A 40417-2022'
$ws.Cells.Item(15, 2).Value = 44823
$ws.Cells.Item(15, 3).Value = 46079
$ws.Cells.Item(15, 6).ClearContents()
$ws.Cells.Item(15, 7).Value = 2.3
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = 0
$ws.Cells.Item(15, 14).Value = 0
$ws.Cells.Item(15, 15).Value = 0
$ws.Cells.Item(15, 16).Value = 0
$ws.Cells.Item(15, 17).Value = 0
$ws.Cells.Item(15, 18).ClearContents()
$ws.Cells.Item(15, 19).ClearContents()
$ws.Cells.Item(15, 20).ClearContents()
$ws.Cells.Item(15, 22).ClearContents()
$ws.Cells.Item(15, 23).ClearContents()
$ws.Cells.Item(15, 24).ClearContents()
$ws.Cells.Item(15, 25).ClearContents()

# Row 16 (source data from original row 22)
$ws.Cells.Item(16, 1).Value = 'A 22195-2023'
$ws.Cells.Item(16, 2).Value = 45069.74605324074
$ws.Cells.Item(16, 3).Value = 46079
$ws.Cells.Item(16, 6).Value = 'Övriga Aktiebolag'
$ws.Cells.Item(16, 7).Value = 1.1
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = 0
$ws.Cells.Item(16, 14).Value = 0
$ws.Cells.Item(16, 15).Value = 0
$ws.Cells.Item(16, 16).Value = 0
$ws.Cells.Item(16, 17).Value = 0
$ws.Cells.Item(16, 18).ClearContents()
$ws.Cells.Item(16, 19).ClearContents()
$ws.Cells.Item(16, 20).ClearContents()
$ws.Cells.Item(16, 22).ClearContents()
$ws.Cells.Item(16, 23).ClearContents()
$ws.Cells.Item(16, 24).ClearContents()
$ws.Cells.Item(16, 25).ClearContents()

# Row 17 (source data from original row 23)
$ws.Cells.Item(17, 1).Value = 'A 34302-2024'
$ws.Cells.Item(17, 2).Value = 45524
$ws.Cells.Item(17, 3).Value = 46079
$ws.Cells.Item(17, 6).ClearContents()
$ws.Cells.Item(17, 7).Value = 1.9
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 0
$ws.Cells.Item(17, 13).Value = 0
$ws.Cells.Item(17, 14).Value = 0
$ws.Cells.Item(17, 15).Value = 0
$ws.Cells.Item(17, 16).Value = 0
$ws.Cells.Item(17, 17).Value = 0
$ws.Cells.Item(17, 18).ClearContents()
$ws.Cells.Item(17, 19).ClearContents()
$ws.Cells.Item(17, 20).ClearContents()
$ws.Cells.Item(17, 22).ClearContents()
$ws.Cells.Item(17, 23).ClearContents()
$ws.Cells.Item(17, 24).ClearContents()
$ws.Cells.Item(17, 25).ClearContents()

# Row 18 (source data from original row 20)
$ws.Cells.Item(18, 1).Value = 'A 38631-2023'
$ws.Cells.Item(18, 2).Value = 45162
$ws.Cells.Item(18, 3).Value = 46079
$ws.Cells.Item(18, 6).ClearContents()
$ws.Cells.Item(18, 7).Value = 0.8
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = 0
$ws.Cells.Item(18, 14).Value = 0
$ws.Cells.Item(18, 15).Value = 0
$ws.Cells.Item(18, 16).Value = 0
$ws.Cells.Item(18, 17).Value = 0
$ws.Cells.Item(18, 18).ClearContents()
$ws.Cells.Item(18, 19).ClearContents()
$ws.Cells.Item(18, 20).ClearContents()
$ws.Cells.Item(18, 22).ClearContents()
$ws.Cells.Item(18, 23).ClearContents()
$ws.Cells.Item(18, 24).ClearContents()
$ws.Cells.Item(18, 25).ClearContents()

# Row 19 (source data from original row 18)
$ws.Cells.Item(19, 1).Value = 'A 18090-2022'
$ws.Cells.Item(19, 2).Value = 44684
$ws.Cells.Item(19, 3).Value = 46079
$ws.Cells.Item(19, 6).ClearContents()
$ws.Cells.Item(19, 7).Value = 4.9
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = 0
$ws.Cells.Item(19, 14).Value = 0
$ws.Cells.Item(19, 15).Value = 0
$ws.Cells.Item(19, 16).Value = 0
$ws.Cells.Item(19, 17).Value = 0
$ws.Cells.Item(19, 18).ClearContents()
$ws.Cells.Item(19, 19).ClearContents()
$ws.Cells.Item(19, 20).ClearContents()
$ws.Cells.Item(19, 22).ClearContents()
$ws.Cells.Item(19, 23).ClearContents()
$ws.Cells.Item(19, 24).ClearContents()
$ws.Cells.Item(19, 25).ClearContents()

# Row 20 (source data from original row 14)
$ws.Cells.Item(20, 1).Value = 'A 49536-2025'
$ws.Cells.Item(20, 2).Value = 45939.4221875
$ws.Cells.Item(20, 3).Value = 46079
$ws.Cells.Item(20, 6).ClearContents()
$ws.Cells.Item(20, 7).Value = 1.5
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = 0
$ws.Cells.Item(20, 14).Value = 0
$ws.Cells.Item(20, 15).Value = 0
$ws.Cells.Item(20, 16).Value = 0
$ws.Cells.Item(20, 17).Value = 0
$ws.Cells.Item(20, 18).ClearContents()
$ws.Cells.Item(20, 19).ClearContents()
$ws.Cells.Item(20, 20).ClearContents()
$ws.Cells.Item(20, 22).ClearContents()
$ws.Cells.Item(20, 23).ClearContents()
$ws.Cells.Item(20, 24).ClearContents()
$ws.Cells.Item(20, 25).ClearContents()

# Row 21 (source data from original row 15)
$ws.Cells.Item(21, 1).Value = 'A 49543-2025'
$ws.Cells.Item(21, 2).Value = 45939.42862268518
$ws.Cells.Item(21, 3).Value = 46079
$ws.Cells.Item(21, 6).ClearContents()
$ws.Cells.Item(21, 7).Value = 1.4
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 13).Value = 0
$ws.Cells.Item(21, 14).Value = 0
$ws.Cells.Item(21, 15).Value = 0
$ws.Cells.Item(21, 16).Value = 0
$ws.Cells.Item(21, 17).Value = 0
$ws.Cells.Item(21, 18).ClearContents()
$ws.Cells.Item(21, 19).ClearContents()
$ws.Cells.Item(21, 20).ClearContents()
$ws.Cells.Item(21, 22).ClearContents()
$ws.Cells.Item(21, 23).ClearContents()
$ws.Cells.Item(21, 24).ClearContents()
$ws.Cells.Item(21, 25).ClearContents()

# Row 22 (source data from original row 16)
$ws.Cells.Item(22, 1).Value = 'A 32596-2024'
$ws.Cells.Item(22, 2).Value = 45513.61667824074
$ws.Cells.Item(22, 3).Value = 46079
$ws.Cells.Item(22, 6).ClearContents()
$ws.Cells.Item(22, 7).Value = 2.6
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = 0
$ws.Cells.Item(22, 14).Value = 0
$ws.Cells.Item(22, 15).Value = 0
$ws.Cells.Item(22, 16).Value = 0
$ws.Cells.Item(22, 17).Value = 0
$ws.Cells.Item(22, 18).ClearContents()
$ws.Cells.Item(22, 19).ClearContents()
$ws.Cells.Item(22, 20).ClearContents()
$ws.Cells.Item(22, 22).ClearContents()
$ws.Cells.Item(22, 23).ClearContents()
$ws.Cells.Item(22, 24).ClearContents()
$ws.Cells.Item(22, 25).ClearContents()

# Row 23 (source data from original row 21)
$ws.Cells.Item(23, 1).Value = 'A 49549-2025'
$ws.Cells.Item(23, 2).Value = 45939
$ws.Cells.Item(23, 3).Value = 46079
$ws.Cells.Item(23, 6).ClearContents()
$ws.Cells.Item(23, 7).Value = 0.5
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 13).Value = 0
$ws.Cells.Item(23, 14).Value = 0
$ws.Cells.Item(23, 15).Value = 0
$ws.Cells.Item(23, 16).Value = 0
$ws.Cells.Item(23, 17).Value = 0
$ws.Cells.Item(23, 18).ClearContents()
$ws.Cells.Item(23, 19).ClearContents()
$ws.Cells.Item(23, 20).ClearContents()
$ws.Cells.Item(23, 22).ClearContents()
$ws.Cells.Item(23, 23).ClearContents()
$ws.Cells.Item(23, 24).ClearContents()
$ws.Cells.Item(23, 25).ClearContents()
